$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update section name 2LCRCP -> 2LCRCP-HOU
$ws.Range("B2").Value = "2LCRCP-HOU"

# Update the active selection to match the saved workbook state
$ws.Range("E5").Select()
